# DOMA-4452: add "Meter place" column for import (Kitchen/Bathroom values)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing header/body cells (column R) onto
# the new column S so the new column matches the look of its neighbours.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$ws.Range("R2:R11").Copy()
$ws.Range("S2:S11").PasteSpecial(-4122)

# Match column width to column R (23.5 in the saved file).
$ws.Range("S1").ColumnWidth = $ws.Range("R1").ColumnWidth

# Header
$ws.Range("S1").Value = "Meter place"

# Body values - alternating Kitchen / Bathroom per row
$ws.Range("S2").Value = "Kitchen"
$ws.Range("S3").Value = "Bathroom"
$ws.Range("S4").Value = "Kitchen"
$ws.Range("S5").Value = "Bathroom"
$ws.Range("S6").Value = "Kitchen"
$ws.Range("S7").Value = "Bathroom"
$ws.Range("S8").Value = "Kitchen"
$ws.Range("S9").Value = "Bathroom"
$ws.Range("S10").Value = "Kitchen"
$ws.Range("S11").Value = "Bathroom"
